$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column B (GoalPoseX/Y and GoalOrient X/Y/Z/W)
$ws.Range("B12").Value = 37.39
$ws.Range("B13").Value = 69.49
$ws.Range("B15").Value = 0.637
$ws.Range("B16").Value = 0.327
$ws.Range("B17").Value = 0.327
$ws.Range("B18").Value = 0.627

# Widen column A slightly (closest the COM layer's rounding allows to 37.6)
$ws.Columns("A").ColumnWidth = 36.8

# Move the active selection to B19
$ws.Range("B19").Select()

# Shrink the sheet-tab area of the window (tabRatio 983 -> 500 in the XML)
$excel.ActiveWindow.TabRatio = 500
